$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.277791050588093
$ws.Range("D2").Value = 0.2294187233413965
$ws.Range("E2").Value = 0.2504832339759027
$ws.Range("F2").Value = 1.477857911800015
$ws.Range("G2").Value = 0.002449362423981744
$ws.Range("J2").Value = 0.3308574286622985
$ws.Range("L2").Value = 0.3048730014661629
$ws.Range("M2").Value = 0.320277874890003
$ws.Range("N2").Value = 1.854056015280008
$ws.Range("O2").Value = 3.607955291981369

$ws.Range("B3").Value = 1.231078397773757
$ws.Range("D3").Value = 0.2294926430559698
$ws.Range("E3").Value = 0.2487000765755418
$ws.Range("F3").Value = 1.478008731710467
$ws.Range("G3").Value = 0.00245269181224168
$ws.Range("J3").Value = 0.326332659866182
$ws.Range("L3").Value = 0.2787639667401862
$ws.Range("M3").Value = 0.3028187283928432
$ws.Range("N3").Value = 1.860525384613396
$ws.Range("O3").Value = 3.585789942055442

$ws.Range("B4").Value = 1.202804816984496
$ws.Range("D4").Value = 0.229561043579448
$ws.Range("E4").Value = 0.2476218287037426
$ws.Range("F4").Value = 1.478956922603338
$ws.Range("G4").Value = 0.002454847845837984
$ws.Range("J4").Value = 0.3235682811492779
$ws.Range("L4").Value = 0.2627520603419384
$ws.Range("M4").Value = 0.2921738771505815
$ws.Range("N4").Value = 1.865177161379066
$ws.Range("O4").Value = 3.574437087433637

$ws.Range("B5").Value = 1.191386723882232
$ws.Range("D5").Value = 0.229594729831593
$ws.Range("E5").Value = 0.2471866833957357
$ws.Range("F5").Value = 1.479558405925651
$ws.Range("G5").Value = 0.002455754639666295
$ws.Range("J5").Value = 0.3224453980770718
$ws.Range("L5").Value = 0.256232337685617
$ws.Range("M5").Value = 0.2878552307457767
$ws.Range("N5").Value = 1.867244117459428
$ws.Range("O5").Value = 3.570377813953542

$ws.Range("B6").Value = 1.189497045579714
$ws.Range("D6").Value = 0.2296006752208015
$ws.Range("E6").Value = 0.2471146865872775
$ws.Range("F6").Value = 1.479671270933252
$ws.Range("G6").Value = 0.002455906917604683
$ws.Range("J6").Value = 0.3222591672293831
$ws.Range("L6").Value = 0.2551500734601859
$ws.Range("M6").Value = 0.2871392929925989
$ws.Range("N6").Value = 1.867597694229786
$ws.Range("O6").Value = 3.569738015729371

$ws.Range("B7").Value = 1.202650407731539
$ws.Range("D7").Value = 0.2295614743142522
$ws.Range("E7").Value = 0.2476159428815983
$ws.Range("F7").Value = 1.478964163612744
$ws.Range("G7").Value = 0.002454859960902218
$ws.Range("J7").Value = 0.3235531226946051
$ws.Range("L7").Value = 0.2626641111977648
$ws.Range("M7").Value = 0.292115556137496
$ws.Range("N7").Value = 1.865204342820618
$ws.Range("O7").Value = 3.574380046961267

$ws.Range("B8").Value = 1.261600382070725
$ws.Range("D8").Value = 0.2294394495811716
$ws.Range("E8").Value = 0.2498649939759083
$ws.Range("F8").Value = 1.477732299063987
$ws.Range("G8").Value = 0.002450487250014621
$ws.Range("J8").Value = 0.3292944976780134
$ws.Range("L8").Value = 0.2958668906911015
$ws.Range("M8").Value = 0.3142425548969427
$ws.Range("N8").Value = 1.856145871010455
$ws.Range("O8").Value = 3.599844045975232

$ws.Range("B9").Value = 1.380403390024441
$ws.Range("D9").Value = 0.2293816863194813
$ws.Range("E9").Value = 0.2544043765283241
$ws.Range("F9").Value = 1.482109933214744
$ws.Range("G9").Value = 0.002442795344279138
$ws.Range("J9").Value = 0.3406575440912363
$ws.Range("L9").Value = 0.3611130683476063
$ws.Range("M9").Value = 0.3582178016231126
$ws.Range("N9").Value = 1.843756052973376
$ws.Range("O9").Value = 3.667710197710733

$ws.Range("B10").Value = 1.469601087480441
$ws.Range("D10").Value = 0.2294485358393601
$ws.Range("E10").Value = 0.2578146335845872
$ws.Range("F10").Value = 1.489476958424262
$ws.Range("G10").Value = 0.002437676899234479
$ws.Range("J10").Value = 0.3490627681119207
$ws.Range("L10").Value = 0.4091150060211817
$ws.Range("M10").Value = 0.3908700307144102
$ws.Range("N10").Value = 1.837906484938273
$ws.Range("O10").Value = 3.728549096719632

$ws.Range("B11").Value = 1.510587200653163
$ws.Range("D11").Value = 0.2295023936507903
$ws.Range("E11").Value = 0.2593816595475147
$ws.Range("F11").Value = 1.493731939622222
$ws.Range("G11").Value = 0.002435462912645148
$ws.Range("J11").Value = 0.3528974780366525
$ws.Range("L11").Value = 0.4309632447030936
$ws.Range("M11").Value = 0.4057965412202265
$ws.Range("N11").Value = 1.83594717148172
$ws.Range("O11").Value = 3.758620678027341

$ws.Range("B12").Value = 1.526165626217733
$ws.Range("D12").Value = 0.2295261355085252
$ws.Range("E12").Value = 0.2599772396323452
$ws.Range("F12").Value = 1.495473266207881
$ws.Range("G12").Value = 0.002434640898434953
$ws.Range("J12").Value = 0.3543510504292726
$ws.Range("L12").Value = 0.4392379326511104
$ws.Range("M12").Value = 0.4114590175959094
$ws.Range("N12").Value = 1.835305744784677
$ws.Range("O12").Value = 3.770353132234504

$ws.Range("B13").Value = 1.5228079759558
$ws.Range("D13").Value = 0.2295208738037999
$ws.Range("E13").Value = 0.2598488749611363
$ws.Range("F13").Value = 1.495092454293129
$ws.Range("G13").Value = 0.002434817206939866
$ws.Range("J13").Value = 0.3540379354535474
$ws.Range("L13").Value = 0.4374557832782671
$ws.Range("M13").Value = 0.4102390569003944
$ws.Range("N13").Value = 1.835439422673545
$ws.Range("O13").Value = 3.767810986323354

$ws.Range("B14").Value = 1.511867693112322
$ws.Range("D14").Value = 0.2295042799681042
$ws.Range("E14").Value = 0.2594306149837706
$ws.Range("F14").Value = 1.493872592438223
$ws.Range("G14").Value = 0.002435394957455548
$ws.Range("J14").Value = 0.3530170360889144
$ws.Range("L14").Value = 0.4316439862677726
$ws.Range("M14").Value = 0.4062621954938237
$ws.Range("N14").Value = 1.835892388437486
$ws.Range("O14").Value = 3.759578996906555

$ws.Range("B15").Value = 1.505173960340471
$ws.Range("D15").Value = 0.2294945509237678
$ws.Range("E15").Value = 0.2591747005616796
$ws.Range("F15").Value = 1.493142332391216
$ws.Range("G15").Value = 0.00243575097622882
$ws.Range("J15").Value = 0.3523918902736796
$ws.Range("L15").Value = 0.4280842393029332
$ws.Range("M15").Value = 0.4038275621702709
$ws.Range("N15").Value = 1.836182922969158
$ws.Range("O15").Value = 3.754581609754496

$ws.Range("B16").Value = 1.466930703285243
$ws.Range("D16").Value = 0.2294454858180259
$ws.Range("E16").Value = 0.2577125340168394
$ws.Range("F16").Value = 1.489217079140218
$ws.Range("G16").Value = 0.00243782388439306
$ws.Range("J16").Value = 0.3488123719831862
$ws.Range("L16").Value = 0.4076873710119457
$ws.Range("M16").Value = 0.3898959845494829
$ws.Range("N16").Value = 1.838048617836719
$ws.Range("O16").Value = 3.726632098235427

$ws.Range("B17").Value = 1.443573886362287
$ws.Range("D17").Value = 0.2294213738522224
$ws.Range("E17").Value = 0.2568195063665399
$ws.Range("F17").Value = 1.487040596113701
$ws.Range("G17").Value = 0.002439124795265297
$ws.Range("J17").Value = 0.3466191977916111
$ws.Range("L17").Value = 0.3951772634533768
$ws.Range("M17").Value = 0.3813678254535731
$ws.Range("N17").Value = 1.839372607743911
$ws.Range("O17").Value = 3.710099924547904

$ws.Range("B18").Value = 1.430178281068208
$ws.Range("D18").Value = 0.2294097137534017
$ws.Range("E18").Value = 0.2563073410321905
$ws.Range("F18").Value = 1.48587379318559
$ws.Range("G18").Value = 0.002439883819441681
$ws.Range("J18").Value = 0.3453587986622608
$ws.Range("L18").Value = 0.3879829306656859
$ws.Range("M18").Value = 0.3764695312101622
$ws.Range("N18").Value = 1.840200199735762
$ws.Range("O18").Value = 3.700816516207396

$ws.Range("B19").Value = 1.425649423427956
$ws.Range("D19").Value = 0.2294061459048926
$ws.Range("E19").Value = 0.2561341871891258
$ws.Range("F19").Value = 1.485493338308174
$ws.Range("G19").Value = 0.002440142664883469
$ws.Range("J19").Value = 0.344932235224114
$ws.Range("L19").Value = 0.3855472653812342
$ws.Range("M19").Value = 0.3748122446456037
$ws.Range("N19").Value = 1.840491767066908
$ws.Range("O19").Value = 3.697712021341374

$ws.Range("B20").Value = 1.446056270048871
$ws.Range("D20").Value = 0.2294237122237739
$ws.Range("E20").Value = 0.25691441798919
$ws.Range("F20").Value = 1.487263483037083
$ws.Range("G20").Value = 0.002438985196666802
$ws.Range("J20").Value = 0.3468525566545821
$ws.Range("L20").Value = 0.3965088701958166
$ws.Range("M20").Value = 0.3822749532715761
$ws.Range("N20").Value = 1.839224831729766
$ws.Range("O20").Value = 3.71183646464732

$ws.Range("B21").Value = 1.515079556895216
$ws.Range("D21").Value = 0.2295090633364367
$ws.Range("E21").Value = 0.2595534094634573
$ws.Range("F21").Value = 1.494227364842757
$ws.Range("G21").Value = 0.002435224814378677
$ws.Range("J21").Value = 0.3533168607912245
$ws.Range("L21").Value = 0.4333510222495534
$ws.Range("M21").Value = 0.4074300238841388
$ws.Range("N21").Value = 1.835756616234193
$ws.Range("O21").Value = 3.761987564540732

$ws.Range("B22").Value = 1.560527122756127
$ws.Range("D22").Value = 0.2295843432172084
$ws.Range("E22").Value = 0.2612908235117324
$ws.Range("F22").Value = 1.499536789636963
$ws.Range("G22").Value = 0.00243286259396641
$ws.Range("J22").Value = 0.3575500454438583
$ws.Range("L22").Value = 0.4574364704004097
$ws.Range("M22").Value = 0.4239292074965348
$ws.Range("N22").Value = 1.834075690064466
$ws.Range("O22").Value = 3.796775360180447

$ws.Range("B23").Value = 1.536240418841089
$ws.Range("D23").Value = 0.2295423889595121
$ws.Range("E23").Value = 0.2603623969892048
$ws.Range("F23").Value = 1.496633645210764
$ws.Range("G23").Value = 0.002434114651487641
$ws.Range("J23").Value = 0.3552899969467376
$ws.Range("L23").Value = 0.4445811335193355
$ws.Range("M23").Value = 0.4151180077849119
$ws.Range("N23").Value = 1.834919358761681
$ws.Range("O23").Value = 3.778024271301376

$ws.Range("B24").Value = 1.44493388294228
$ws.Range("D24").Value = 0.2294226481856221
$ws.Range("E24").Value = 0.256871504552219
$ws.Range("F24").Value = 1.487162452694434
$ws.Range("G24").Value = 0.002439048274495159
$ws.Range("J24").Value = 0.3467470535890271
$ws.Range("L24").Value = 0.3959068573095408
$ws.Range("M24").Value = 0.3818648262488082
$ws.Range("N24").Value = 1.839291434364014
$ws.Range("O24").Value = 3.711050685999538

$ws.Range("B25").Value = 1.347925252702538
$ws.Range("D25").Value = 0.2293780219623507
$ws.Range("E25").Value = 0.2531629024101001
$ws.Range("F25").Value = 1.480197426433776
$ws.Range("G25").Value = 0.00244478224975839
$ws.Range("J25").Value = 0.337573079700654
$ws.Range("L25").Value = 0.3434495032098255
$ws.Range("M25").Value = 0.3462600856176365
$ws.Range("N25").Value = 1.846534907831369
$ws.Range("O25").Value = 3.647426157327487
